$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 596
$ws.Range("I31").Value = 620
$ws.Range("K31").Value = 1860
$ws.Range("M31").Value = -1630

$ws.Range("H51").Value = 2381.9092
$ws.Range("I51").Value = 2360.2
$ws.Range("J51").Value = 2400
$ws.Range("K51").Value = 2360.2
$ws.Range("L51").Value = 2400
$ws.Range("M51").Value = -1876.2
$ws.Range("N51").Value = -3368

$ws.Range("H125").Value = 1182.75
$ws.Range("J125").Value = 1270.3334
$ws.Range("L125").Value = 11433.0006
$ws.Range("N125").Value = -16353.0006

$ws.Range("H137").Value = 13097978
$ws.Range("I137").Value = 5191
$ws.Range("J137").Value = 20371748
$ws.Range("K137").Value = 15573
$ws.Range("L137").Value = 61115244
$ws.Range("M137").Value = -13023
$ws.Range("N137").Value = -61120344


# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 5000
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -5230

$ws.Range("H32").Value = 604992.4
$ws.Range("I32").Value = 29009.158
$ws.Range("J32").Value = 2168375.5
$ws.Range("K32").Value = 29009.158
$ws.Range("L32").Value = 2168375.5
$ws.Range("M32").Value = -28722.158
$ws.Range("N32").Value = -2168949.5

$ws.Range("H61").Value = 21746158
$ws.Range("I61").Value = 35722070
$ws.Range("K61").Value = 35722070
$ws.Range("M61").Value = -35721858

$ws.Range("H110").Value = 2067.7778
$ws.Range("I110").Value = 2067.7778
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 2067.7778
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -22.77779999999984
$ws.Range("N110").ClearContents()

$ws.Range("H120").Value = 45800
$ws.Range("J120").Value = 45800
$ws.Range("L120").Value = 45800
$ws.Range("N120").Value = -55476

$ws.Range("H132").Value = 3757.2727
$ws.Range("I132").Value = 3869.0715
$ws.Range("J132").Value = 3561.625
$ws.Range("K132").Value = 11607.2145
$ws.Range("L132").Value = 10684.875
$ws.Range("M132").Value = -9077.2145
$ws.Range("N132").Value = -15744.875

$ws.Range("H136").Value = 21746158
$ws.Range("I136").Value = 35722070
$ws.Range("K136").Value = 107166210
$ws.Range("M136").Value = -107163660


# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 226.6842
$ws.Range("I80").Value = 117.166664
$ws.Range("J80").Value = 277.23077
$ws.Range("K80").Value = 117.166664
$ws.Range("L80").Value = 277.23077
$ws.Range("M80").Value = 880.833336
$ws.Range("N80").Value = -2273.23077

$ws.Range("H83").Value = 226.6842
$ws.Range("I83").Value = 117.166664
$ws.Range("J83").Value = 277.23077
$ws.Range("K83").Value = 585.83332
$ws.Range("L83").Value = 1386.15385
$ws.Range("M83").Value = 4406.16668
$ws.Range("N83").Value = -11370.15385

$ws.Range("H134").Value = 804.5714
$ws.Range("I134").Value = 819.8
$ws.Range("J134").Value = 500
$ws.Range("K134").Value = 2459.4
$ws.Range("L134").Value = 1500
$ws.Range("M134").Value = 75.60000000000036
$ws.Range("N134").Value = -6570


# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 2000
$ws.Range("I36").Value = 2000
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2000
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1612
$ws.Range("N36").ClearContents()

$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1840
$ws.Range("N40").ClearContents()

$ws.Range("H58").Value = 2397041.5
$ws.Range("I58").Value = 3789564.5
$ws.Range("J58").Value = 9859.643
$ws.Range("K58").Value = 3789564.5
$ws.Range("L58").Value = 9859.643
$ws.Range("M58").Value = -3789361.5
$ws.Range("N58").Value = -10265.643

$ws.Range("H99").Value = 1314.2858
$ws.Range("I99").Value = 1316.6666
$ws.Range("J99").Value = 1300
$ws.Range("K99").Value = 1316.6666
$ws.Range("L99").Value = 1300
$ws.Range("M99").Value = 181.3334
$ws.Range("N99").Value = -4296

$ws.Range("H126").Value = 1314.2858
$ws.Range("I126").Value = 1316.6666
$ws.Range("J126").Value = 1300
$ws.Range("K126").Value = 3949.9998
$ws.Range("L126").Value = 3900
$ws.Range("M126").Value = -1479.9998
$ws.Range("N126").Value = -8840

$ws.Range("H132").Value = 4694.5264
$ws.Range("I132").Value = 4546.533
$ws.Range("K132").Value = 13639.599
$ws.Range("M132").Value = -11109.599

$ws.Range("H136").Value = 2397041.5
$ws.Range("I136").Value = 3789564.5
$ws.Range("J136").Value = 9859.643
$ws.Range("K136").Value = 11368693.5
$ws.Range("L136").Value = 29578.929
$ws.Range("M136").Value = -11366143.5
$ws.Range("N136").Value = -34678.929


# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2635
$ws.Range("I17").Value = 270
$ws.Range("K17").Value = 810
$ws.Range("M17").Value = -641

$ws.Range("H61").Value = 258.26666
$ws.Range("I61").Value = 46.9
$ws.Range("J61").Value = 681
$ws.Range("K61").Value = 140.7
$ws.Range("L61").Value = 2043
$ws.Range("M61").Value = 74.30000000000001
$ws.Range("N61").Value = -2473

$ws.Range("H68").Value = 3142.8914
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 3142.8914
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 9428.674199999999
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -11050.6742

$ws.Range("H71").Value = 3142.8914
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 3142.8914
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 28286.0226
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -36398.0226

$ws.Range("H101").Value = 7007.25
$ws.Range("J101").Value = 7007.25
$ws.Range("L101").Value = 21021.75
$ws.Range("N101").Value = -25889.75

$ws.Range("H113").Value = 419.3402
$ws.Range("I113").Value = 489.85715
$ws.Range("J113").Value = 365.4909
$ws.Range("K113").Value = 1469.57145
$ws.Range("L113").Value = 1096.4727
$ws.Range("M113").Value = 700.4285500000001
$ws.Range("N113").Value = -5436.4727

$ws.Range("H131").Value = 1069.0605
$ws.Range("I131").Value = 247.71428
$ws.Range("J131").Value = 1290.1923
$ws.Range("K131").Value = 743.14284
$ws.Range("L131").Value = 3870.5769
$ws.Range("M131").Value = 4296.85716
$ws.Range("N131").Value = -13950.5769

$ws.Range("H132").Value = 2611.6216
$ws.Range("J132").Value = 2406.389
$ws.Range("L132").Value = 21657.501
$ws.Range("N132").Value = -26717.501


# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3795.7
$ws.Range("I102").Value = 3158.3333
$ws.Range("J102").Value = 4751.75
$ws.Range("K102").Value = 3158.3333
$ws.Range("L102").Value = 4751.75
$ws.Range("M102").Value = -1536.3333
$ws.Range("N102").Value = -7995.75

$ws.Range("H137").Value = 45000
$ws.Range("J137").Value = 45000
$ws.Range("L137").Value = 45000
$ws.Range("N137").Value = -55200


# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 14706.235
$ws.Range("I61").Value = 15968
$ws.Range("K61").Value = 15968
$ws.Range("M61").Value = -15766

$ws.Range("H113").Value = 14706.235
$ws.Range("I113").Value = 15968
$ws.Range("K113").Value = 15968
$ws.Range("M113").Value = -13798

$ws.Range("H136").Value = 3487.4736
$ws.Range("I136").Value = 2583.1177
$ws.Range("K136").Value = 7749.353099999999
$ws.Range("M136").Value = -5199.353099999999


# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 21400000
$ws.Range("I21").Value = 21400000
$ws.Range("K21").Value = 21400000
$ws.Range("M21").Value = -21399765

$ws.Range("H28").Value = 30900
$ws.Range("J28").Value = 30900
$ws.Range("L28").Value = 30900
$ws.Range("N28").Value = -31596

$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H35").Value = 21400000
$ws.Range("I35").Value = 21400000
$ws.Range("K35").Value = 21400000
$ws.Range("M35").Value = -21399710

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H43").Value = 14999.5
$ws.Range("J43").Value = 14999.5
$ws.Range("L43").Value = 14999.5
$ws.Range("N43").Value = -15297.5

$ws.Range("H48").Value = 13065
$ws.Range("J48").Value = 13065
$ws.Range("L48").Value = 13065
$ws.Range("N48").Value = -14203

$ws.Range("H49").Value = 20062
$ws.Range("J49").Value = 20062
$ws.Range("L49").Value = 20062
$ws.Range("N49").Value = -20522

$ws.Range("H50").Value = 10067
$ws.Range("I50").Value = 10067
$ws.Range("K50").Value = 10067
$ws.Range("M50").Value = -9436

$ws.Range("H51").Value = 8666.666999999999
$ws.Range("I51").Value = 8666.666999999999
$ws.Range("K51").Value = 8666.666999999999
$ws.Range("M51").Value = -8156.666999999999

$ws.Range("H113").Value = 638.7826
$ws.Range("I113").Value = 260.41666
$ws.Range("J113").Value = 1051.5454
$ws.Range("K113").Value = 781.2499799999999
$ws.Range("L113").Value = 3154.6362
$ws.Range("M113").Value = 1388.75002
$ws.Range("N113").Value = -7494.6362

$ws.Range("H132").Value = 3416.4119
$ws.Range("I132").Value = 3327.0715
$ws.Range("J132").Value = 3833.3333
$ws.Range("K132").Value = 9981.2145
$ws.Range("L132").Value = 11499.9999
$ws.Range("M132").Value = -7451.2145
$ws.Range("N132").Value = -16559.9999

$ws.Range("H136").Value = 5041.146
$ws.Range("I136").Value = 1948.2273
$ws.Range("J136").Value = 7658.231
$ws.Range("K136").Value = 5844.6819
$ws.Range("L136").Value = 22974.693
$ws.Range("M136").Value = -3294.6819
$ws.Range("N136").Value = -28074.693

